# Update scripts with new TPM values.
# The underlying source data changed from "ECs"/"FAPs" sending-cluster rows
# (7 data rows) down to just the "FAPs" sending-cluster rows (3 data rows),
# with refreshed TPM-derived statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: FAPs -> Cd28/Cd80 -> ECs -----------------------------------
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cd28"
$ws.Range("C2").Value = "Cd80"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.042039
$ws.Range("H2").Value = 0.126117
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.08106466666666667
$ws.Range("N2").Value = 0.243194
$ws.Range("O2").Value = 0.01252465659474717
$ws.Range("P2").Value = 0.01252465659474717
$ws.Range("Q2").Value = 0.003407877522
$ws.Range("R2").Value = 0.030670897698
$ws.Range("S2").Value = 0.01252465659474717
$ws.Range("T2").Value = 0.01252465659474717

# --- Row 3: FAPs -> Cd28/Cd80 -> FAPs -----------------------------------
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cd28"
$ws.Range("C3").Value = "Cd80"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.042039
$ws.Range("H3").Value = 0.126117
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.993807666666666
$ws.Range("N3").Value = 17.981423
$ws.Range("O3").Value = 0.926055528343168
$ws.Range("P3").Value = 0.926055528343168
$ws.Range("Q3").Value = 0.251973680499
$ws.Range("R3").Value = 2.267763124491
$ws.Range("S3").Value = 0.926055528343168
$ws.Range("T3").Value = 0.926055528343168

# --- Row 4: FAPs -> Cd28/Cd80 -> MuSCs ----------------------------------
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cd28"
$ws.Range("C4").Value = "Cd80"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.042039
$ws.Range("H4").Value = 0.126117
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.397534
$ws.Range("N4").Value = 1.192602
$ws.Range("O4").Value = 0.06141981506208485
$ws.Range("P4").Value = 0.06141981506208484
$ws.Range("Q4").Value = 0.016711931826
$ws.Range("R4").Value = 0.150407386434
$ws.Range("S4").Value = 0.06141981506208485
$ws.Range("T4").Value = 0.06141981506208484

# --- Remove the old rows 5:7 (FAPs rows that previously lived there) ---
$ws.Rows("5:7").Delete()
